{"js": "// Apply the Bescheid.docx revision described in the commit diff.\n// Strategy: locate each changed sentence/paragraph with Range.search()\n// (matching on exact old text, including the \\u000b \"soft line break\"\n// characters produced by <w:br/> when a whole paragraph is removed) and\n// replace it in place with Range.insertText(..., Word.InsertLocation.replace).\n\nconst body = context.document.body;\nconst VT = \"\\u000b\"; // text representation of <w:br/>\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${oldText.substring(0, 40)}...\" but found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Tenor item 1: reworded sentence.\nawait replaceOnce(\n  \"1. Sie sind verpflichtet, das Fachwerkhaus in Neuried, Lange Stra\u00dfe 12, mit Biberschwanz-Dachziegeln zu reparieren.\",\n  \"1. Sie sind als Eigent\u00fcmer des Fachwerkhauses in Neuried, Lange Stra\u00dfe 12, verpflichtet, das Dach mit Biberschwanz-Dachziegeln zu reparieren.\"\n);\n\n// 2. Begr\u00fcndung intro paragraph: reworded + extra sentence about legal basis.\nawait replaceOnce(\n  \"Das Fachwerkhaus ist ein Kulturdenkmal im Sinne von \u00a7 2 Abs. 1 DSchG, da dessen Erhaltung aus heimatgeschichtlichen Gr\u00fcnden ein \u00f6ffentliches Interesse besteht. Durch den Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt, wodurch eine Gef\u00e4hrdung des Denkmals vorliegt.\",\n  \"Das Fachwerkhaus stellt ein Kulturdenkmal dar, da es aus heimatgeschichtlichen Gr\u00fcnden ein \u00f6ffentliches Interesse an seiner Erhaltung gibt (\u00a7 2 Abs. 1 DSchG). Durch den Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt, wodurch eine Gef\u00e4hrdung des Denkmals vorliegt. Die Reparaturanordnung st\u00fctzt sich auf \u00a7 1 Abs. 1 in Verbindung mit \u00a7 7 Abs. 1 Satz 1 DSchG sowie \u00a7 7 PolG.\"\n);\n\n// 3. Remove the whole \"Als Pflichtige kommen...\" paragraph, together with\n//    the pair of line breaks that separated it from the preceding paragraph\n//    (the pair of breaks that follows it remains, now separating the two\n//    neighbouring paragraphs).\nawait replaceOnce(\n  VT + VT +\n    \"Als Pflichtige kommen sowohl Sie als auch Ihr Bruder Georg Konrad in Betracht. Sie sind nach \u00a7 7 Abs. 1 Satz 1 DSchG und \u00a7 7 PolG als Eigent\u00fcmer des Fachwerkhauses verpflichtet, von dem eine Gefahr ausgeht. Ihr Bruder Georg Konrad ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\",\n  \"\"\n);\n\n// 4. Replace the \"Die Anordnung der BSD...\" paragraph with the new merged\n//    paragraph (new wording + content moved from the removed paragraph #3).\nawait replaceOnce(\n  \"Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen. Der Vorteil f\u00fcr die Allgemeinheit durch die Ansehnlichkeit des Denkmals rechtfertigt den finanziellen Nachteil f\u00fcr den Eigent\u00fcmer.\",\n  \"Als Eigent\u00fcmer des Fachwerkhauses sind Sie verpflichtet, die Reparatur durchzuf\u00fchren. Eine kosteng\u00fcnstigere Reparatur mit Eternitplatten ist nicht geeignet, die Denkmalanforderungen zu erf\u00fcllen. Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da der Vorteil f\u00fcr die Allgemeinheit durch die Ansehnlichkeit des Denkmals den finanziellen Nachteil f\u00fcr den Eigent\u00fcmer \u00fcberwiegt.\"\n);\n\n// 5. \"Eine privatrechtliche Unm\u00f6glichkeit...\" paragraph: reworded + extra sentence.\nawait replaceOnce(\n  \"Eine privatrechtliche Unm\u00f6glichkeit liegt nicht vor, da Sie ohne die Mitwirkung Ihres Bruders handeln k\u00f6nnen und die Ma\u00dfnahme als notwendige Erhaltungsma\u00dfnahme anzusehen ist.\",\n  \"Eine privatrechtliche Unm\u00f6glichkeit liegt nicht vor, da Sie ohne die Mitwirkung Ihres Bruders Georg Konrad handeln k\u00f6nnen. Eine Mitbestimmungspflicht nach \u00a7 2038 Abs. 1 Satz 1 BGB entf\u00e4llt, da die Reparaturanordnung als notwendige Erhaltungsma\u00dfnahme anzusehen ist (\u00a7 2038 Abs. 1 Satz 2 Halbsatz 2 BGB).\"\n);\n\n// 6. Remove the whole \"Die Anordnung ist bestimmt genug...\" paragraph, along\n//    with the pair of line breaks that precede it.\nawait replaceOnce(\n  VT + VT +\n    \"Die Anordnung ist bestimmt genug formuliert und entspricht den formellen Voraussetzungen.\",\n  \"\"\n);\n\n// 7. Rechtsbehelfsbelehrung paragraph: reworded references.\nawait replaceOnce(\n  \"Gegen die Dachdeckungsanordnung k\u00f6nnen Sie innerhalb eines Monats nach Bekanntgabe Widerspruch einlegen (\u00a7 37 Abs. 6 LVwVfG, \u00a7 70 VwGO). Gegen die Anordnung der sofortigen Vollziehung k\u00f6nnen Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht stellen (\u00a7 80 Abs. 5 VwGO).\",\n  \"Gegen die Reparaturanordnung k\u00f6nnen Sie innerhalb eines Monats nach Bekanntgabe Widerspruch einlegen (\u00a7 37 Abs. 6 LVwVfG). Gegen die Anordnung der sofortigen Vollziehung k\u00f6nnen Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht Freiburg stellen (\u00a7 80 Abs. 5 VwGO).\"\n);\n", "ps1": "# Apply the Bescheid.docx revision described in the commit diff using the\n# Word COM object model. $word / $d (ActiveDocument) are pre-seeded.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Sentence([string]$oldText, [string]$newText) {\n    # Find-and-replace a single exact occurrence of $oldText with $newText.\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Replace-Sentence: text not found: $oldText\"\n    }\n}\n\nfunction Remove-ParagraphWithPrecedingBreaks([string]$searchText) {\n    # Delete $searchText together with the two manual line breaks (<w:br/>)\n    # immediately preceding it, collapsing the paragraph out of the body.\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Remove-ParagraphWithPrecedingBreaks: text not found: $searchText\"\n    }\n    $range.MoveStart(1, -2)\n    $range.Text = \"\"\n}\n\n# 1. Tenor item 1: reworded sentence.\nReplace-Sentence `\n    \"1. Sie sind verpflichtet, das Fachwerkhaus in Neuried, Lange Stra\u00dfe 12, mit Biberschwanz-Dachziegeln zu reparieren.\" `\n    \"1. Sie sind als Eigent\u00fcmer des Fachwerkhauses in Neuried, Lange Stra\u00dfe 12, verpflichtet, das Dach mit Biberschwanz-Dachziegeln zu reparieren.\"\n\n# 2. Begr\u00fcndung intro paragraph: reworded + extra sentence about legal basis.\nReplace-Sentence `\n    \"Das Fachwerkhaus ist ein Kulturdenkmal im Sinne von \u00a7 2 Abs. 1 DSchG, da dessen Erhaltung aus heimatgeschichtlichen Gr\u00fcnden ein \u00f6ffentliches Interesse besteht. Durch den Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt, wodurch eine Gef\u00e4hrdung des Denkmals vorliegt.\" `\n    \"Das Fachwerkhaus stellt ein Kulturdenkmal dar, da es aus heimatgeschichtlichen Gr\u00fcnden ein \u00f6ffentliches Interesse an seiner Erhaltung gibt (\u00a7 2 Abs. 1 DSchG). Durch den Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt, wodurch eine Gef\u00e4hrdung des Denkmals vorliegt. Die Reparaturanordnung st\u00fctzt sich auf \u00a7 1 Abs. 1 in Verbindung mit \u00a7 7 Abs. 1 Satz 1 DSchG sowie \u00a7 7 PolG.\"\n\n# 3. Remove the whole \"Als Pflichtige kommen...\" paragraph, together with the\n#    pair of line breaks separating it from the preceding paragraph.\nRemove-ParagraphWithPrecedingBreaks `\n    \"Als Pflichtige kommen sowohl Sie als auch Ihr Bruder Georg Konrad in Betracht. Sie sind nach \u00a7 7 Abs. 1 Satz 1 DSchG und \u00a7 7 PolG als Eigent\u00fcmer des Fachwerkhauses verpflichtet, von dem eine Gefahr ausgeht. Ihr Bruder Georg Konrad ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\"\n\n# 4. Replace the \"Die Anordnung der BSD...\" paragraph with the new merged\n#    paragraph (new wording + content moved from the removed paragraph #3).\nReplace-Sentence `\n    \"Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen. Der Vorteil f\u00fcr die Allgemeinheit durch die Ansehnlichkeit des Denkmals rechtfertigt den finanziellen Nachteil f\u00fcr den Eigent\u00fcmer.\" `\n    \"Als Eigent\u00fcmer des Fachwerkhauses sind Sie verpflichtet, die Reparatur durchzuf\u00fchren. Eine kosteng\u00fcnstigere Reparatur mit Eternitplatten ist nicht geeignet, die Denkmalanforderungen zu erf\u00fcllen. Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da der Vorteil f\u00fcr die Allgemeinheit durch die Ansehnlichkeit des Denkmals den finanziellen Nachteil f\u00fcr den Eigent\u00fcmer \u00fcberwiegt.\"\n\n# 5. \"Eine privatrechtliche Unm\u00f6glichkeit...\" paragraph: reworded + extra sentence.\nReplace-Sentence `\n    \"Eine privatrechtliche Unm\u00f6glichkeit liegt nicht vor, da Sie ohne die Mitwirkung Ihres Bruders handeln k\u00f6nnen und die Ma\u00dfnahme als notwendige Erhaltungsma\u00dfnahme anzusehen ist.\" `\n    \"Eine privatrechtliche Unm\u00f6glichkeit liegt nicht vor, da Sie ohne die Mitwirkung Ihres Bruders Georg Konrad handeln k\u00f6nnen. Eine Mitbestimmungspflicht nach \u00a7 2038 Abs. 1 Satz 1 BGB entf\u00e4llt, da die Reparaturanordnung als notwendige Erhaltungsma\u00dfnahme anzusehen ist (\u00a7 2038 Abs. 1 Satz 2 Halbsatz 2 BGB).\"\n\n# 6. Remove the whole \"Die Anordnung ist bestimmt genug...\" paragraph, along\n#    with the pair of line breaks that precede it.\nRemove-ParagraphWithPrecedingBreaks `\n    \"Die Anordnung ist bestimmt genug formuliert und entspricht den formellen Voraussetzungen.\"\n\n# 7. Rechtsbehelfsbelehrung paragraph: reworded references.\nReplace-Sentence `\n    \"Gegen die Dachdeckungsanordnung k\u00f6nnen Sie innerhalb eines Monats nach Bekanntgabe Widerspruch einlegen (\u00a7 37 Abs. 6 LVwVfG, \u00a7 70 VwGO). Gegen die Anordnung der sofortigen Vollziehung k\u00f6nnen Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht stellen (\u00a7 80 Abs. 5 VwGO).\" `\n    \"Gegen die Reparaturanordnung k\u00f6nnen Sie innerhalb eines Monats nach Bekanntgabe Widerspruch einlegen (\u00a7 37 Abs. 6 LVwVfG). Gegen die Anordnung der sofortigen Vollziehung k\u00f6nnen Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht Freiburg stellen (\u00a7 80 Abs. 5 VwGO).\"\n\nWrite-Output \"DONE\"\n"}
